$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 140: new form submission
$ws.Cells.Item(140, 1).Value = "1219丶X"
$ws.Cells.Item(140, 2).NumberFormat = "yyyy/m/d h:mm:ss;@"
$ws.Cells.Item(140, 2).Value = 46018.6778587963
$ws.Cells.Item(140, 3).Value = "12a6b740"
$ws.Cells.Item(140, 4).NumberFormat = "@"
$ws.Cells.Item(140, 4).Value = "1833704413"
$ws.Cells.Item(140, 4).ClearFormats()

# Row 141: new form submission
$ws.Cells.Item(141, 1).Value = "不因如此"
$ws.Cells.Item(141, 2).NumberFormat = "yyyy/m/d h:mm:ss;@"
$ws.Cells.Item(141, 2).Value = 46018.8404050926
$ws.Cells.Item(141, 3).Value = "887ae519"
$ws.Cells.Item(141, 4).NumberFormat = "@"
$ws.Cells.Item(141, 4).Value = "2399377548"
$ws.Cells.Item(141, 4).ClearFormats()
